# Fruta / hortaliza, semanal
# Insert a new weekly price record as a new row 147 (pushing existing
# rows 147:160 down to 148:161), duplicating the former row 147's
# attributes but with an updated date (2023-07-25 / serial 45132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 147, shifting rows 147:160
# down to 148:161.
$ws.Rows.Item(147).Insert()

# Populate the new row 147 with the same record as the old row 147
# (Agrícola del Norte S.A. de Arica / Naranja / Fukumoto / Segunda /
# Región de O'Higgins) but with the new report date.
$ws.Range("A147").Value = 1
$ws.Range("B147").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C147").Value = "Arica y Parinacota"
$ws.Range("D147").Value = 45132
$ws.Range("D147").NumberFormat = $ws.Range("D148").NumberFormat
$ws.Range("E147").Value = 15
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100102
$ws.Range("H147").Value = "Cítricos"
$ws.Range("I147").Value = 100102005
$ws.Range("J147").Value = "Naranja"
$ws.Range("K147").Value = "Fukumoto"
$ws.Range("L147").Value = "Segunda"
$ws.Range("M147").Value = 270
$ws.Range("N147").Value = 900
$ws.Range("O147").Value = 950
$ws.Range("P147").Value = 925
$ws.Range("Q147").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R147").Value = "Región de O'Higgins"
$ws.Range("S147").Value = 925
$ws.Range("T147").Value = 1
